$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a 5th week column (E) mirroring the existing weekly date columns (B:D).
# E3 gets the next weekly date (7 days after D3); copy D3's format (style) first
# so the same cellXf (date number format) is reused instead of minting a new one.
$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value2 = $ws.Range("D3").Value2 + 7

# E4:E8 get the same "attended" marker (Wingdings "ü" checkmark) and style as D4:D8.
$ws.Range("D4:D8").Copy()
$ws.Range("E4:E8").PasteSpecial(-4122)
$ws.Range("E4:E8").Value2 = $ws.Range("D4").Value2

# Match the selection recorded in the saved view state.
$ws.Range("E4:E8").Select()
